$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2178217821782178
$ws.Range("C2").Value = 0.5082508250825083
$ws.Range("J2").Value = 0.0132013201320132
$ws.Range("P2").Value = 0.1848184818481848
$ws.Range("S2").Value = 0.07590759075907591

# Row 3
$ws.Range("B3").Value = 0.01910828025477707
$ws.Range("C3").Value = 0.006369426751592357
$ws.Range("J3").Value = 0.01910828025477707
$ws.Range("P3").Value = 0.7452229299363057
$ws.Range("S3").Value = 0.2101910828025478

# Row 4
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.723404255319149
$ws.Range("S4").Value = 0.2553191489361702

# Row 6
$ws.Range("B6").Value = 0.09677419354838709
$ws.Range("D6").Value = 0.02580645161290323
$ws.Range("F6").Value = 0.04516129032258064
$ws.Range("J6").Value = 0.232258064516129
$ws.Range("O6").Value = 0.01290322580645161
$ws.Range("Q6").Value = 0.1548387096774194
$ws.Range("R6").Value = 0.07741935483870968
$ws.Range("S6").Value = 0.3548387096774194

# Row 7
$ws.Range("B7").Value = 0.0898876404494382
$ws.Range("D7").Value = 0.02808988764044944
$ws.Range("F7").Value = 0.0449438202247191
$ws.Range("J7").Value = 0.1348314606741573
$ws.Range("O7").Value = 0.005617977528089887
$ws.Range("Q7").Value = 0.1573033707865168
$ws.Range("R7").Value = 0.07865168539325842
$ws.Range("S7").Value = 0.4606741573033708

# Row 8
$ws.Range("B8").Value = 0.1037234042553191
$ws.Range("D8").Value = 0.01595744680851064
$ws.Range("F8").Value = 0.04521276595744681
$ws.Range("J8").Value = 0.1515957446808511
$ws.Range("O8").Value = 0.01861702127659574
$ws.Range("Q8").Value = 0.148936170212766
$ws.Range("R8").Value = 0.1063829787234043
$ws.Range("S8").Value = 0.4095744680851064

# Row 9
$ws.Range("B9").Value = 0.09895833333333333
$ws.Range("D9").Value = 0.02083333333333333
$ws.Range("F9").Value = 0.05729166666666666
$ws.Range("J9").Value = 0.09895833333333333
$ws.Range("O9").Value = 0.03125
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.1041666666666667
$ws.Range("S9").Value = 0.421875

# Row 10
$ws.Range("B10").Value = 0.1357913669064748
$ws.Range("D10").Value = 0.02607913669064748
$ws.Range("F10").Value = 0.05485611510791367
$ws.Range("J10").Value = 0.1088129496402878
$ws.Range("O10").Value = 0.01528776978417266
$ws.Range("Q10").Value = 0.2068345323741007
$ws.Range("R10").Value = 0.09172661870503597
$ws.Range("S10").Value = 0.3606115107913669

# Row 11
$ws.Range("G11").Value = 0.1097560975609756
$ws.Range("J11").Value = 0.1056910569105691
$ws.Range("K11").Value = 0.1747967479674797
$ws.Range("L11").Value = 0.5894308943089431
$ws.Range("S11").Value = 0.02032520325203252

# Row 12
$ws.Range("G12").Value = 0.7804878048780488
$ws.Range("J12").Value = 0.1585365853658537
$ws.Range("L12").Value = 0.03658536585365853
$ws.Range("S12").Value = 0.02439024390243903

# Row 13
$ws.Range("G13").Value = 0.7105263157894737
$ws.Range("J13").Value = 0.1842105263157895
$ws.Range("S13").Value = 0.1052631578947368

# Row 15
$ws.Range("F15").Value = 0.01036269430051814
$ws.Range("H15").Value = 0.1295336787564767
$ws.Range("I15").Value = 0.09326424870466321
$ws.Range("J15").Value = 0.4093264248704663
$ws.Range("K15").Value = 0.05181347150259067
$ws.Range("M15").Value = 0.005181347150259068
$ws.Range("O15").Value = 0.07772020725388601
$ws.Range("S15").Value = 0.2227979274611399

# Row 16
$ws.Range("F16").Value = 0.009852216748768473
$ws.Range("H16").Value = 0.1083743842364532
$ws.Range("I16").Value = 0.06896551724137931
$ws.Range("J16").Value = 0.4975369458128079
$ws.Range("K16").Value = 0.06896551724137931
$ws.Range("M16").Value = 0.01970443349753695
$ws.Range("N16").Value = 0.004926108374384237
$ws.Range("O16").Value = 0.09852216748768473
$ws.Range("S16").Value = 0.1231527093596059

# Row 17
$ws.Range("F17").Value = 0.005420054200542005
$ws.Range("H17").Value = 0.1707317073170732
$ws.Range("I17").Value = 0.07859078590785908
$ws.Range("J17").Value = 0.4254742547425474
$ws.Range("K17").Value = 0.1084010840108401
$ws.Range("M17").Value = 0.02168021680216802
$ws.Range("O17").Value = 0.05962059620596206
$ws.Range("S17").Value = 0.1300813008130081

# Row 18
$ws.Range("F18").Value = 0.0160427807486631
$ws.Range("H18").Value = 0.1657754010695187
$ws.Range("I18").Value = 0.08021390374331551
$ws.Range("J18").Value = 0.3796791443850268
$ws.Range("K18").Value = 0.0855614973262032
$ws.Range("M18").Value = 0.0213903743315508
$ws.Range("N18").Value = 0.0053475935828877
$ws.Range("O18").Value = 0.09090909090909091
$ws.Range("S18").Value = 0.1550802139037433

# Row 19
$ws.Range("F19").Value = 0.02001740644038294
$ws.Range("H19").Value = 0.1940818102697998
$ws.Range("I19").Value = 0.1018276762402089
$ws.Range("J19").Value = 0.3594429939077459
$ws.Range("K19").Value = 0.1053089643167972
$ws.Range("M19").Value = 0.01827676240208877
$ws.Range("O19").Value = 0.0557006092254134
$ws.Range("S19").Value = 0.1453437771975631
